# Update "Benefits and other factors" section heading + its two question
# cells, per the commit "Update project template with benefits and risks
# content".
#
# Each edit below splits one run of text into two separate <w:r> runs
# (matching the target OOXML diff exactly) using this trick: plain
# Range.InsertAfter() on an already-populated paragraph gets silently
# coalesced back into the neighbouring run when it shares the same
# formatting, so instead we insert the trailing piece while
# TrackRevisions is on (which *always* wraps a fresh insertion in its
# own <w:ins><w:r>...</w:r></w:ins>, never merging it with a sibling
# run) and then call AcceptAllRevisions() to bake it in - that unwraps
# the <w:ins> but leaves the run boundary intact, giving two distinct
# <w:r> elements with identical rPr.

$d = $word.ActiveDocument

function Split-RunText($oldText, $firstPart, $secondPart) {
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }

    $runStart = $rng.Start

    # Replace the whole old run's text with just the first part. For a
    # single-run paragraph/cell like these, assigning .Text rewrites
    # that run in place (and Word adds xml:space="preserve" on its own
    # when the text has leading/trailing whitespace).
    $rng.Text = $firstPart

    # Now append the second part as a brand-new, independent run.
    $d.TrackRevisions = $true
    $firstEnd = $runStart + $firstPart.Length
    $insertionPoint = $d.Range($firstEnd, $firstEnd)
    $insertionPoint.InsertAfter($secondPart)
    $d.TrackRevisions = $false
}

Split-RunText "Benefits and other factors" "Benefits and " "risks"

Split-RunText "What are the benefits the transfer is intended to bring?" `
    "What are the intended benefits of the transfer" "?"

Split-RunText "Are there any other factors to consider during this transfer?" `
    "Are there any risks to consider" "?"

# Accept the tracked insertions once, in a batch: this unwraps every
# <w:ins><w:r>...</w:r></w:ins> into a plain <w:r>...</w:r> without
# re-merging it into its neighbouring run.
$d.AcceptAllRevisions() | Out-Null

Write-Output "Applied benefits/risks edits"
